$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46065
$ws.Range("D2").Value = -0.24
$ws.Range("G2").Value = -0.22
$ws.Range("I2").Value = 1.46
$ws.Range("J2").Value = 2.97
$ws.Range("K2").Value = 0.45
$ws.Range("L2").Value = 0.01
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.01
$ws.Range("S2").Value = 0.57
$ws.Range("T2").Value = 10.81
$ws.Range("U2").Value = 32.62
$ws.Range("V2").Value = 35
$ws.Range("W2").Value = 29.15
$ws.Range("X2").Value = 18.51
$ws.Range("Y2").Value = 7.95
$ws.Range("Z2").Value = 5.76
$ws.Range("AB2").Value = 22.65
$ws.Range("AD2").Value = 32.08
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 21.72
$ws.Range("AG2").Value = "0h-17h"
